$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "flashcards" to "Sheet"
$ws.Name = "Sheet"

# Append a new flashcard entry in row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Conhecimentos Específicos"
$ws.Cells.Item(2, 3).Value = "Gestão da Qualidade"
$ws.Cells.Item(2, 4).Value = "processos universais da Trilogia de Juran"
$ws.Cells.Item(2, 5).Value = "<ul>`n<li>Planejamento da Qualidade</li>`n<li>Controle da Qualidade</li>`n<li>Melhoria da Qualidade</li>`n</ul>"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
